$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the H1 title.
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Cut()

# ---------------------------------------------------------------------------
# 2) Paste that paragraph back in, right before the final ("Create a feature
#    image...") paragraph, so it becomes a new paragraph just above it.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertionPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insertionPoint.Paste()

# Cut/Paste drops the paragraph's leading empty run (<w:r/>), so restore it
# to keep the paragraph shape consistent with the rest of the document.
$newCount = $d.Paragraphs.Count
$pastedPara = $d.Paragraphs.Item($newCount - 1)
$leadRange = $d.Range($pastedPara.Range.Start, $pastedPara.Range.Start)
$emptyRunXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$leadRange.InsertXML($emptyRunXml)

# ---------------------------------------------------------------------------
# 3) Turn that pasted paragraph into the new bold
#    "Play Baron Samedi Free: A Unique Voodoo Slot Game" line (replacing its
#    "Meta description: ..." text, but keeping the bold run formatting).
# ---------------------------------------------------------------------------
$pastedPara = $d.Paragraphs.Item($newCount - 1)
$boldTextRange = $d.Range($pastedPara.Range.Start, $pastedPara.Range.End - 1)
$boldTextRange.Text = "Play Baron Samedi Free: A Unique Voodoo Slot Game"

# ---------------------------------------------------------------------------
# 4) Replace the text of the final (italic, "Create a feature image...")
#    paragraph with the meta-description copy.
# ---------------------------------------------------------------------------
$imagePromptPara = $d.Paragraphs.Item($newCount)
$imagePromptRange = $d.Range($imagePromptPara.Range.Start, $imagePromptPara.Range.End - 1)
$imagePromptRange.Text = "Experience the mysticism of voodoo with Baron Samedi slot game, featuring collectible cards and high RTP. Play for free now."
